$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the D1 header text (title above the data table).
$ws.Range("D1").Value = "Average Percent MSA Residues Masked (Average Number of Residues Masked)"

# 2. Update the footnote in A14 to the new text, with "number" in italics.
$ws.Range("A14").Value = " Percent gaps were calculated from unfiltered amino-acid alignments as the total number of gaps divided by the total number of MSA positions, and represent the percentage of columns with at least one gap, averaged across all MSA replicates. As these values were calculated from amino-acid alignments, the number of sites masked should be multiplied by 3 for codon alignments, but the percentage stays the same."

$text = $ws.Range("A14").Text
$start = $text.IndexOf("number of sites masked") + 1
$len = "number".Length
$ws.Range("A14").Characters($start, $len).Font.Italic = $true

# 3. Move the selection to A15 (cosmetic change matching the diff).
$ws.Range("A15").Select()
